$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new year column (P) to the table, matching the formatting of
# the adjacent column O (year 2021) by copying its format into column P.
$ws.Range("O3:O5").Copy()
$ws.Range("P3:P5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New data for 2022
$ws.Range("P3").Value = 2022
$ws.Range("P4").Value = 15
$ws.Range("P5").Value = 2130.4

# Move the active selection to P6, as in the saved workbook
$null = $ws.Range("P6").Select()
